$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "@prefix" sheet: add three missing prefix declarations (pixels, channel,
#    bindata) that mirror the existing "image" prefix row, and set the
#    sheet's page setup (paper size / orientation).
# ---------------------------------------------------------------------------
$wsPrefix = $wb.Worksheets.Item("@prefix")

$wsPrefix.Range("A11").Value = "pixels"
$wsPrefix.Range("B11").Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/samples/time-series/pixels/"

$wsPrefix.Range("A12").Value = "channel"
$wsPrefix.Range("B12").Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/samples/time-series/channel/"

$wsPrefix.Range("A13").Value = "bindata"
$wsPrefix.Range("B13").Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/samples/time-series/bindata/"

$wsPrefix.PageSetup.PaperSize = 9
$wsPrefix.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 2. "Image" sheet: the sample pixels reference now uses the namespaced
#    "pixels:" naming convention instead of the old "[pixels:0:0]" form.
# ---------------------------------------------------------------------------
$wsImage = $wb.Worksheets.Item("Image")
$wsImage.Range("E5").Value = "pixels:pixels0.0"

# ---------------------------------------------------------------------------
# 3. "Pixels" sheet: sample rows now use namespaced identifiers for pixels,
#    channel and bindata references.
# ---------------------------------------------------------------------------
$wsPixels = $wb.Worksheets.Item("Pixels")
$wsPixels.Range("B5").Value = "pixels:pixels0.0"
$wsPixels.Range("B6").Value = "pixels:pixels0.0"
$wsPixels.Range("B7").Value = "pixels:pixels0.0"
$wsPixels.Range("B8").Value = "pixels:pixels0.0"
$wsPixels.Range("B9").Value = "pixels:pixels0.0"

$wsPixels.Range("K5").Value = "channel:channel1"

$wsPixels.Range("L5").Value = "bindata:bindata1"
$wsPixels.Range("L6").Value = "bindata:bindata2"
$wsPixels.Range("L7").Value = "bindata:bindata3"
$wsPixels.Range("L8").Value = "bindata:bindata4"
$wsPixels.Range("L9").Value = "bindata:bindata5"

# ---------------------------------------------------------------------------
# 4. "Channel" sheet: sample channel reference renamed the same way.
# ---------------------------------------------------------------------------
$wsChannel = $wb.Worksheets.Item("Channel")
$wsChannel.Range("B5").Value = "channel:channel1"

# ---------------------------------------------------------------------------
# 5. "Binary_Data" sheet: sample bindata references renamed the same way.
# ---------------------------------------------------------------------------
$wsBinData = $wb.Worksheets.Item("Binary_Data")
$wsBinData.Range("B5").Value = "bindata:bindata1"
$wsBinData.Range("B6").Value = "bindata:bindata2"
$wsBinData.Range("B7").Value = "bindata:bindata3"
$wsBinData.Range("B8").Value = "bindata:bindata4"
$wsBinData.Range("B9").Value = "bindata:bindata5"
